# Adds "Errors" and "Warnings" sheets after "Classes", with an error
# message written into the Errors sheet, and makes the Errors sheet the
# active tab (matching what Excel persists when a user adds sheets and
# leaves the new one selected).

$wb = $excel.ActiveWorkbook
$classes = $wb.Worksheets.Item("Classes")

# Insert "Errors" right after "Classes", then "Warnings" right after "Errors".
$errors = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $classes)
$errors.Name = "Errors"

$warnings = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $errors)
$warnings.Name = "Warnings"

# Populate the Errors sheet with the missing-class message.
$errors.Range("A1").Value = 'Sheet "Classes" Row: 2 Missing "OFF CLS"'

# Leave "Errors" as the active/selected sheet.
$errors.Activate()
